# Applies the 'Book Synthesis' edit: turns the single paragraph into a
# centered, underlined 'Summary' heading followed by six body paragraphs
# summarizing 'Between the World and Me' by Ta-Nehisi Coates.
$d = $word.ActiveDocument

# --- Paragraph 1: replace the old text in place with 'Summary' ---
$p1 = $d.Paragraphs(1)
$r = $d.Range($p1.Range.Start, $p1.Range.End)
$r.Text = 'Summary'

# --- Paragraph 2: split off a new paragraph and fill its runs ---
$pPrev = $d.Paragraphs(1)
$rEnd = $d.Range($pPrev.Range.Start, $pPrev.Range.End)
$rEnd.InsertParagraphAfter()
$pNew = $d.Paragraphs(2)
$cursor = $d.Range($pNew.Range.Start, $pNew.Range.Start)
$cursor.InsertAfter('Ta-Nehisi Coates wrote the book Between the World and Me as a letter to his son, Samori, explaining his experiences growing up ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('i')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('n a white America that was built on the backs of blacks with violence and oppression. The major roots of American culture were explained to be tied to enslavement, injustice, theft and violence. The beginning of the book is a recount of his time on a talk show where the host ask him about his views on American progress. Ta-Nehisi ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('explains that ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('America has subjugated groups of people into their natural born races and complexions rather than their national backgrounds which has led European decent Americans into believing they are white over their nationality and ultimately adopt the American white culture of violence amongst those who are not alike. For when a country has a justice system that preys on and hunts a group of Americans for the color of their skin with no repercussion, how can the progress of America not be built on violence and racism. ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('He explains how the black body has been created fragile in the American progression where systemic oppression and racism aims to hurt this body. ')
$cursor = $d.Range($cursor.End, $cursor.End)

# --- Paragraph 3: split off a new paragraph and fill its runs ---
$pPrev = $d.Paragraphs(2)
$rEnd = $d.Range($pPrev.Range.Start, $pPrev.Range.End)
$rEnd.InsertParagraphAfter()
$pNew = $d.Paragraphs(3)
$cursor = $d.Range($pNew.Range.Start, $pNew.Range.Start)
$cursor.InsertAfter('Growing up in an African American community in Baltimore, Ta-Nehisi witnessed how systematic oppression shaped the everyday lives of those around him. ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('Fear consumed the everyday customs of the people in his community, and it was witnessed ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('from their garments')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter(' to their music, and to their parenting styles.')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter(' ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('After comparing his living community to that of the white community that he had seen on TV and experienced in person, he realized that his part of the world felt like a whole other galaxy. Before discovering the world and venturing out into the public, Ta-Nehisi had to survive the streets of his community where few make it out unscathed.')
$cursor = $d.Range($cursor.End, $cursor.End)

# --- Paragraph 4: split off a new paragraph and fill its runs ---
$pPrev = $d.Paragraphs(3)
$rEnd = $d.Range($pPrev.Range.Start, $pPrev.Range.End)
$rEnd.InsertParagraphAfter()
$pNew = $d.Paragraphs(4)
$cursor = $d.Range($pNew.Range.Start, $pNew.Range.Start)
$cursor.InsertAfter('The 1990s is where Ta-Nehisi found his inspiration to do good as a writer with the inspiration of Malcolm X within black culture where his resurgence in pop culture and media exploded. ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('From this deep found idolization of one of his ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('communities’')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter(' strongest leaders, Ta-Nehisi hoped to one ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('day accomplish his own study, exploration and publishment of his ideas through books. Howard University was the college that Ta-Nehisi chose to continue his higher education where he felt empowered and passionate with the like-minded culture of intelligent and impactful African Americans just like him. This is where he realized the world was more than just the white ')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('America')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter(' he grew up in. All his time was spent learning and reading about his ancestral culture and upbringing where he discovered and unmasked all the history that was smoke screened by his boyhood education and community.')
$cursor = $d.Range($cursor.End, $cursor.End)

# --- Paragraph 5: split off a new paragraph and fill its runs ---
$pPrev = $d.Paragraphs(4)
$rEnd = $d.Range($pPrev.Range.Start, $pPrev.Range.End)
$rEnd.InsertParagraphAfter()
$pNew = $d.Paragraphs(5)
$cursor = $d.Range($pNew.Range.Start, $pNew.Range.Start)
$cursor.InsertAfter('Before Samori was born, Ta-Nehisi recounts the killing of an African American, Prince Jones, that was pulled over by a police officer. Reflecting on the impact of the corrupt justice system and the effects of the Prince Jones killing, Ta-Nehisi chooses to respond with his words through writing.')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter(' This event worries him on the injustices that his son will have to grow up into as a black American where the reality of "the Dream" for white Americans comes at the cost of black Americans. ')
$cursor = $d.Range($cursor.End, $cursor.End)

# --- Paragraph 6: split off a new paragraph and fill its runs ---
$pPrev = $d.Paragraphs(5)
$rEnd = $d.Range($pPrev.Range.Start, $pPrev.Range.End)
$rEnd.InsertParagraphAfter()
$pNew = $d.Paragraphs(6)
$cursor = $d.Range($pNew.Range.Start, $pNew.Range.Start)
$cursor.InsertAfter('Ta-Nehisi teaches his son Samori to remember the generations of black people who were born into chains and to remember it with all its nuance, error, and humanity. Through the history of America and his recounted experience')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter('s, he')
$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.InsertAfter(' teaches his son about their place being built as the essential below of the country but is for Samori to grow up without measuring himself against the white measuring stick.')
$cursor = $d.Range($cursor.End, $cursor.End)

# --- Apply heading formatting (centered + underlined) to paragraph 1 LAST so
#     it does not get inherited by the paragraphs split off above ---
$p1final = $d.Paragraphs(1)
$p1final.Alignment = 1
$p1final.Range.Font.Underline = 1

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
